$wb = $excel.ActiveWorkbook

# --- Sheet "data" (sheet1): add column AG for 24. 8. 2021 ---
$ws1 = $wb.Worksheets.Item("data")

# Copy header formatting from AF1 into AG1, then set header values
$ws1.Range("AF1").Copy() | Out-Null
$ws1.Range("AG1").PasteSpecial(-4122) | Out-Null
$ws1.Range("AG1").Value = "24. 8. 2021"

$data1 = @{
    2 = 0.19
    3 = 0.12
    4 = 0.44
    5 = 0.28
    6 = 0.13
    7 = 0.19
    8 = 0.21
    9 = 0.18
    10 = 0.21
    11 = 0.18
    12 = 0.18
    13 = 0.39
    14 = 0.18
    15 = 0.18
    16 = 0.21
    17 = 0.17
    18 = 0.25
    19 = 0.24
    20 = 0.15
    21 = 0.14
    22 = 0.11
    23 = 0.23
    24 = 0.39
    25 = 0.38
    26 = 0.13
    27 = 0.08
    28 = 0.12
    29 = 0.2
    30 = 0.1
    31 = 0.12
    32 = 0.12
    33 = 0.2
    34 = 0.17
    35 = 0.1
    36 = 0.14
    37 = 0.13
    38 = 0.07000000000000001
    39 = 0.22
    40 = 0.13
    41 = 0.09
    42 = 0.08
    43 = 0.08
    44 = 0.19
    45 = 0.23
}
foreach ($r in $data1.Keys) {
    $ws1.Cells.Item([int]$r, 33).Value = $data1[$r]
}

$ws1.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 1. 9. 2021"

# --- Sheet "pocetR" (sheet2): add column AF for 24. 8. 2021 ---
$ws2 = $wb.Worksheets.Item("pocetR")

$ws2.Range("AE1").Copy() | Out-Null
$ws2.Range("AF1").PasteSpecial(-4122) | Out-Null
$ws2.Range("AF1").Value = "24. 8. 2021"

$data2 = @{
    2 = 1901
    3 = 200
    4 = 365
    5 = 1336
    6 = 916
    7 = 169
    8 = 544
    9 = 272
    10 = 874
    11 = 159
    12 = 112
    13 = 756
    14 = 872
    15 = 654
    16 = 375
    17 = 200
    18 = 729
    19 = 583
    20 = 263
    21 = 598
    22 = 311
    23 = 164
}
foreach ($r in $data2.Keys) {
    $ws2.Cells.Item([int]$r, 32).Value = $data2[$r]
}

# Empty trailing cell AF24 (matches blank inlineStr cells across row 24)
$ws2.Range("AE24").Copy() | Out-Null
$ws2.Range("AF24").PasteSpecial(-4122) | Out-Null

$ws2.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 9. 2021"

Write-Host "edit complete"
